# "Generate Report for Handback"
#
# The handback-status report got regenerated: the zh-cn and de-de sheets
# each record a "Latest Handback DateTime" (column K) for the first data
# row (row 2, the 55b480a7-... file). A fresh handback came in for both
# locales, so those two timestamps advance to the new handback times.
# Everything else in the report (file names, statuses, the other rows,
# etc.) is unchanged.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Latest Handback DateTime for 55b480a7-...zh-cn.xlf
$wsZhCn.Range("K2").Value = "2016-11-08 23:41:20"

# Latest Handback DateTime for 55b480a7-...de-de.xlf
$wsDeDe.Range("K2").Value = "2016-11-08 23:41:37"
